$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Pagos" (F) and "Inscrições homologadas" (H) columns for the
# affected rows. H = F + G (Isenções deferidas stays unchanged).

$updates = @(
    @{ Row = 11; F = 382; H = 447 },
    @{ Row = 12; F = 619; H = 705 },
    @{ Row = 15; F = 131; H = 183 },
    @{ Row = 16; F = 165; H = 213 },
    @{ Row = 27; F = 285; H = 367 },
    @{ Row = 32; F = 173; H = 211 },
    @{ Row = 41; F = 321; H = 413 },
    @{ Row = 48; F = 187; H = 231 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 6).Value = $u.F
    $ws.Cells.Item($u.Row, 8).Value = $u.H
}

$wb.Save()
